$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row: Day, Time, Module Code, Module Title, Hours, Class Type, Lecturer, Room, Block, Group, Level, Course
$ws.Range("H1").Value = "Room"
$ws.Range("I1").Value = "Block"
$ws.Range("J1").Value = "Group"
$ws.Range("K1").Value = "Level"
$ws.Range("L1").Value = "Course"

# Data rows (2-10): A Day, B Time, C Module Code, D Module Title, E Hours, F Class Type, G Lecturer, H Room, I Block, J Group, K Level, L Course
$data = @(
    @("MON", "7:00-9:30",  "5CS024", "Collaborative Development",                 2.5, "Workshop", "Mr. Raj Shrestha",    "Lab-04 Patan",    "HCK", "L5CG5",         5, "BCS"),
    @("TUE", "7:00-9:00",  "5CS022", "Human Computer Interaction",                2,   "Lecture",  "Mr. Apurba Neupane",  "LT-02 Telford",   "WLV", "L5CG(5+6+7+8)", 5, "BCS"),
    @("TUE", "9:30-11:30", "5CS020", "Distributed and Cloud Systems Programming", 2,   "Lecture",  "Mr. Sumanta Silwal",  "LT-01 Wulfruna",  "WLV", "L5CG(5+6+7+8)", 5, "BCS"),
    @("WED", "7:00-9:00",  "5CS024", "Collaborative Development",                 2,   "Lecture",  "Mr. Raj Shrestha",    "LT-02 Telford",   "WLV", "L5CG(5+6+7+8)", 5, "BCS"),
    @("WED", "9:30-11:30", "5CS020", "Distributed and Cloud Systems Programming", 2,   "Tutorial", "Mr. Shishir Poudel",  "TR-01 Dudley",    "WLV", "L5CG5",         5, "BCS"),
    @("THU", "9:00-11:00", "5CS022", "Human Computer Interaction",                2,   "Tutorial", "Mr. Pravash Karki",   "TR-01 Dudley",    "WLV", "L5CG5",         5, "BCS"),
    @("THU", "12:00-14:30","5CS020", "Distributed and Cloud Systems Programming", 2.5, "Workshop", "Mr. Shishir Poudel",  "Lab-02 Moseley",  "WLV", "L5CG5",         5, "BCS"),
    @("FRI", "7:00-9:00",  "5CS024", "Collaborative Development",                 2,   "Tutorial", "Mr. Raj Shrestha",    "TR-01 Dudley",    "WLV", "L5CG5",         5, "BCS"),
    @("FRI", "9:30-12:00", "5CS022", "Human Computer Interaction",                2.5, "Workshop", "Mr. Pravash Karki",   "TR-11 Nagarjung", "HCK", "L5CG5",         5, "BCS")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
    $r++
}
